$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (A = empadronador, B = total_registros)
$data = @(
    @("FERNANDEZ VALDERAS ERNESTO ALI", 136),
    @("GUTIERREZ CARLOS TERESA DE JESUS", 131),
    @("CONTRERAS VALDERRAMA JULIA ALEJANDRA", 120),
    @("VALLE MAGALLAN EDUAR", 119),
    @("ROMERO CHANAME YOSSELY TRINIDAD", 100),
    @("ZAVALETA MANAY JORGE LUIS", 100),
    @("HIDALGO CUBAS LUISA YVONE", 98),
    @("CAMACHO LINARES JUDITH ARLETT", 93),
    @("HUMPIRE CASTILLO IRWIN DEIMER", 92),
    @("SEVERINO AVALOS MARJORIE ISABEL", 91),
    @("SENADOR ARBOLEDA GIANCARLOS EXEBIO", 91),
    @("BALLENA ESQUÉN ASTRID CAROLINA", 89),
    @("ZEVALLOS PACHECO ZOILA XIMENA", 83),
    @("SALAZAR FLORES ANA LIZETH", 10)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
